# Sample Project / Main.xlsx — save edit
# Cell B11 on the "Rules" sheet previously held the text "R40" (row for
# rule R40 / "Good Night"). The commit changes that cell's content to the
# literal text "1" (a new shared string), keeping it as text rather than
# letting Excel auto-convert the digit into a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"      # force text storage so "1" isn't stored as a number
$cell.Value = "1"
